$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Account Advanced Find View")
$ws.Activate()

$ws.Range("D2").Value = "Employer-Simple"

$ws.Range("D4").Select()
